$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CALENDAR")

# --- New header columns: media_path (I), hashtags (J) ---
$ws.Cells.Item(1, 9).Value  = "media_path"
$ws.Cells.Item(1, 10).Value = "hashtags"

# --- New auto-generated draft rows (12-16) ---
$newRows = @(
    @{
        post_id    = "gen_1771098233726_888"
        account_id = "acc_samuel"
        line_id    = "auto_gen"
        date       = "2026-02-15 12:00"
        status     = "draft"
        text       = "[ANALYSIS] Regarding Life: Critical implications emerging."
        media_path = ""
        hashtags   = "#Life"
    },
    @{
        post_id    = "gen_1771098233727_555"
        account_id = "acc_mariate"
        line_id    = "auto_gen"
        date       = "2026-02-15 12:00"
        status     = "draft"
        text       = "Life. That's it. That's the tweet."
        media_path = ""
        hashtags   = "#Life"
    },
    @{
        post_id    = "gen_1771098233727_7"
        account_id = "acc_daniel"
        line_id    = "auto_gen"
        date       = "2026-02-15 12:00"
        status     = "draft"
        text       = "Just thinking about Life... 🤔"
        media_path = ""
        hashtags   = "#Life"
    },
    @{
        post_id    = "gen_1771098233727_279"
        account_id = "acc_nguerrero"
        line_id    = "auto_gen"
        date       = "2026-02-15 12:00"
        status     = "draft"
        text       = "just thinking about life.. 🤔"
        media_path = ""
        hashtags   = "#Life"
    },
    @{
        post_id    = "gen_1771098233727_3"
        account_id = "acc_revistavoces"
        line_id    = "auto_gen"
        date       = "2026-02-15 12:00"
        status     = "draft"
        text       = "BREAKING: Life just changed everything. 🧵👇"
        media_path = ""
        hashtags   = "#Life"
    }
)

$r = 12
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.post_id      # A: post_id
    $ws.Cells.Item($r, 2).Value = $row.account_id   # B: account_id
    $ws.Cells.Item($r, 3).Value = $row.line_id      # C: line_id
    $ws.Cells.Item($r, 4).Value = $row.date         # D: scheduled_date
    $ws.Cells.Item($r, 5).Value = $row.status       # E: status
    $ws.Cells.Item($r, 6).Value = $row.text         # F: content_text
    # G (target_url) and H (action_type) are intentionally left blank for
    # these auto-generated rows.
    $ws.Cells.Item($r, 9).Value  = $row.media_path  # I: media_path (empty string)
    $ws.Cells.Item($r, 10).Value = $row.hashtags    # J: hashtags
    $r++
}

# --- Extend the "number stored as text" ignored-error range to cover the
#     newly added columns/rows (A1:J16), mirroring the sheet's existing
#     text-import convention. ---
$ws.Range("A1:J16").Errors.Item(9).Ignore = $true
